# Updated cryptos list values (Price / Volume(1h)) per source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" "68.819.14"
Set-TextCell "E2" "  -0.22%  "
Set-TextCell "D3" "3.917.23"
Set-TextCell "E4" "  -0.04%  "
Set-TextCell "D5" "604.35"
Set-TextCell "E5" "  +0.38%  "
Set-TextCell "D6" "165.66"
Set-TextCell "E6" "  +0.20%  "
Set-TextCell "D7" "3.915.20"
Set-TextCell "E7" "  +4.51%  "
Set-TextCell "E8" "  +0.13%  "
Set-TextCell "E9" "  -1.44%  "
Set-TextCell "E10" "  -3.70%  "
Set-TextCell "D11" "6.40"
Set-TextCell "E11" "  +0.22%  "
Set-TextCell "D13" "37.26"
Set-TextCell "E13" "  -1.01%  "
Set-TextCell "E14" "  -0.60%  "
Set-TextCell "D15" "4.574.73"
Set-TextCell "E15" "  +4.58%  "
Set-TextCell "D16" "3.926.73"
Set-TextCell "E16" "  +4.47%  "
Set-TextCell "D17" "68.949.31"
Set-TextCell "E17" "  +0.02%  "
Set-TextCell "D18" "7.47"
Set-TextCell "E18" "  +0.94%  "
Set-TextCell "D20" "17.04"
Set-TextCell "E20" "  -3.80%  "
Set-TextCell "D21" "11.13"
Set-TextCell "E21" "  -0.53%  "
Set-TextCell "D22" "487.71"
Set-TextCell "E22" "  -0.42%  "
Set-TextCell "D23" "0.722"
Set-TextCell "E23" "  -0.20%  "
Set-TextCell "E24" "  +12.27%  "
Set-TextCell "D25" "84.41"
Set-TextCell "D26" "2.26"
Set-TextCell "E26" "  -0.22%  "
Set-TextCell "D27" "12.09"
Set-TextCell "E27" "  -1.41%  "
Set-TextCell "D28" "10.15"
Set-TextCell "E28" "  +1.26%  "
Set-TextCell "E29" "  +0.00%  "
Set-TextCell "E30" "  -0.47%  "
Set-TextCell "D31" "4.070.48"
Set-TextCell "E31" "  +4.46%  "
Set-TextCell "E32" "  -1.08%  "
Set-TextCell "D33" "7.85"
Set-TextCell "E33" "  -3.49%  "
Set-TextCell "D34" "32.18"
Set-TextCell "E34" "  +1.96%  "
Set-TextCell "D35" "3.866.07"
Set-TextCell "E35" "  +4.97%  "
Set-TextCell "E36" "  +0.08%  "
Set-TextCell "E37" "  +2.96%  "
Set-TextCell "D38" "0.140"
Set-TextCell "E38" "  +1.59%  "
Set-TextCell "D39" "5.90"
Set-TextCell "E39" "  -0.38%  "
Set-TextCell "E40" "  +0.00%  "
Set-TextCell "D41" "0.321"
Set-TextCell "E41" "  -1.10%  "
Set-TextCell "E42" "  -2.64%  "
Set-TextCell "D43" "436.90"
Set-TextCell "E43" "  +3.25%  "
Set-TextCell "E44" "  +0.61%  "
Set-TextCell "D45" "48.46"
Set-TextCell "E45" "  -0.25%  "
Set-TextCell "D46" "8.49"
Set-TextCell "E46" "  +0.75%  "
Set-TextCell "E47" "  -0.01%  "
Set-TextCell "D48" "2.846.29"
Set-TextCell "E48" "  +2.38%  "
Set-TextCell "D49" "26.24"
Set-TextCell "E49" "  +11.05%  "
Set-TextCell "D50" "141.97"
Set-TextCell "E50" "  +0.19%  "
Set-TextCell "E51" "  +1.43%  "
